$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = [ordered]@{
    "D2" = "66.186.82"
    "E2" = "  +1.39%  "
    "D3" = "3.565.41"
    "E3" = "  +5.36%  "
    "E4" = "  +0.15%  "
    "D5" = "'607.30"
    "E5" = "  +2.50%  "
    "D6" = "'144.85"
    "E6" = "  +2.91%  "
    "D7" = "3.563.82"
    "E7" = "  +5.43%  "
    "E8" = "  +0.24%  "
    "E9" = "  +4.39%  "
    "D10" = "'0.136"
    "E10" = "  +2.79%  "
    "D11" = "'7.93"
    "E11" = "  +0.41%  "
    "E12" = "  +1.76%  "
    "D13" = "4.168.50"
    "E13" = "  +5.62%  "
    "E14" = "  +4.93%  "
    "D15" = "'30.13"
    "E15" = "  +1.96%  "
    "D16" = "3.562.35"
    "E16" = "  +5.51%  "
    "D17" = "66.282.82"
    "E17" = "  +1.53%  "
    "E18" = "  -0.61%  "
    "D19" = "'11.48"
    "E19" = "  +11.10%  "
    "D20" = "'6.21"
    "E20" = "  +2.08%  "
    "D21" = "'14.91"
    "E21" = "  +1.96%  "
    "D22" = "'431.19"
    "E22" = "  +4.28%  "
    "E23" = "  +5.87%  "
    "D24" = "'78.79"
    "E24" = "  +2.22%  "
    "D25" = "3.704.59"
    "E25" = "  +5.40%  "
    "E26" = "  +0.05%  "
    "D27" = "'0.0000118"
    "E27" = "  +10.76%  "
    "E28" = "  +5.16%  "
    "D29" = "'8.04"
    "E29" = "  +4.08%  "
    "D30" = "'9.13"
    "E30" = "  -0.47%  "
    "E31" = "  -0.01%  "
    "E32" = "  +2.60%  "
    "E33" = "  -0.42%  "
    "B34" = "EthereumClassic"
    "C34" = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
    "D34" = "'25.49"
    "E34" = "  +5.21%  "
    "B35" = "RenzoRestakedETH"
    "C35" = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
    "D35" = "3.558.11"
    "E35" = "  +5.44%  "
    "D36" = "'1.76"
    "E36" = "  +5.48%  "
    "E37" = "  +0.06%  "
    "D38" = "'7.91"
    "E38" = "  +5.64%  "
    "D39" = "'5.63"
    "E39" = "  +2.75%  "
    "D40" = "'0.999"
    "E40" = "  +0.19%  "
    "D41" = "'170.06"
    "E41" = "  +3.31%  "
    "D42" = "'0.0856"
    "E42" = "  +0.53%  "
    "D43" = "'5.21"
    "E43" = "  +4.39%  "
    "D44" = "'0.897"
    "E44" = "  +3.56%  "
    "E45" = "  +1.72%  "
    "D46" = "'46.15"
    "E46" = "  +2.08%  "
    "E47" = "  +4.12%  "
    "D48" = "'25.89"
    "E48" = "  -1.82%  "
    "D49" = "'2.37"
    "E49" = "  +5.88%  "
    "D50" = "'7.14"
    "E50" = "  +2.00%  "
    "D51" = "'23.36"
    "E51" = "  +16.33%  "
}

foreach ($ref in $changes.Keys) {
    $ws.Range($ref).Value = $changes[$ref]
}
